$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$predictedQuery = '"spinal cord compression" OR "spinal stenosis" OR (cervical myelopathy) OR "gait disturbance" OR "degenerative disc disease" OR "spinal alignment" OR "myelopathy symptoms" OR "nerve root compression" OR "cervical spondylosis" OR "cervical spondylotic myelopathy" OR "upper motor neuron signs" OR ' + "`n" + '(("mri imaging" OR "cervical spine" OR "neurological examination" OR "emg studies" OR "reflex changes" OR "sensory loss" OR "decompression surgery") AND (Myelopathy))' + "`n"

# Row 14: Predicted / Cervical Myelopathy query results
$ws.Range("A14").Value = "Predicted"
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("B14").Value = $predictedQuery
$ws.Range("C14").Value = 0.872
$ws.Range("D14").Value = 0.334
$ws.Range("E14").Value = 0.66
$ws.Range("F14").Value = 0.47
$ws.Range("G14").Value = 0.745
$ws.Rows(14).EntireRow.AutoFit()

# Row 15: Baseline / Cervical Myelopathy
$ws.Range("A15").Value = "Baseline"
$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("B15").Value = "Cervical Myelopathy"
$ws.Range("C15").Value = 0.787
$ws.Range("D15").Value = 0.636
$ws.Range("E15").Value = 0.752
$ws.Range("F15").Value = 0.484
$ws.Range("G15").Value = 0.699

$excel.CutCopyMode = $false
